# "edit projects, temperature defn, KM's role"
#
# 1. ColumnHeadersNcp!B7 had the wrong wording ("salinity" instead of
#    "temperature") in the attributeDefinition for the temperature row -
#    fix the text.
# 2. Personnel sheet gets a new row for Kate Morkeski as "metadata Provider".

$wb = $excel.ActiveWorkbook

# --- Fix the temperature attribute definition on ColumnHeadersNcp ---
$ncp = $wb.Worksheets.Item("ColumnHeadersNcp")
$ncp.Range("B7").Value = "Underway thermosalinograph temperature in degrees Celsius. URI http://vocab.nerc.ac.uk/collection/P01/current/TEMPSZ01/"
$ncp.Range("B7").Select()

# --- Add Kate Morkeski to the Personnel sheet ---
$personnel = $wb.Worksheets.Item("Personnel")
$personnel.Cells.Item(9, 1).Value  = "Kate"
$personnel.Cells.Item(9, 3).Value  = "Morkeski"
$personnel.Cells.Item(9, 4).Value  = "Northeast U.S. Shelf LTER"
$personnel.Cells.Item(9, 5).Value  = "kmorkeski@whoi.edu"
$personnel.Cells.Item(9, 6).Value  = "0000-0002-2903-5851"
$personnel.Cells.Item(9, 7).Value  = "metadata Provider"
$personnel.Cells.Item(9, 8).Value  = "Northeast U.S. Shelf LTER"
$personnel.Cells.Item(9, 9).Value  = "NSF"
$personnel.Cells.Item(9, 10).Value = "OCE-2322676"

$personnel.Range("A9:J9").Select()
